# Losfahren nach Haltestelle ermoeglichen (in Arbeit)
# Insert a new data row (after the row currently holding 7.3 / 50 / 0 / 0 / 0)
# so the table gains an extra breakpoint: 5.1 km @ 50 km/h.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 currently holds 7.3 / 50 / 0 / 0 / 0 — push it (and everything below)
# down by one row, then fill the freed row with the new 5.1 entry.
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value = 5.0999999999999996
$ws.Cells.Item(6, 2).Value = 50
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0

# Copy the style (border formatting, no fill) used by the rest of the table
# into the newly inserted row.
$ws.Range("A7:E7").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection, matching the saved workbook state.
$ws.Range("G8").Select()
